# plot.xlsx — "Slight changes in table demo columns"
#
# The "survey" sheet's choice_list_name column (D5) pointed at the
# choice-list "plants"; rename it to "planting" (the choices sheet still
# uses the per-choice values corn/bean/cotton - those are untouched).
#
# Also restore the view state left behind when the edit was made: the
# "survey" sheet (with the cell that was just edited) becomes the active
# sheet/tab with D5 selected, while the previously-active "settings" sheet
# keeps a plain (unselected) cell reference where the cursor was left.

$wb = $excel.ActiveWorkbook

$survey = $wb.Worksheets.Item("survey")
$settings = $wb.Worksheets.Item("settings")

# Core content edit: choice_list_name "plants" -> "planting"
$survey.Range("D5").Value = "planting"

# View-state: "survey" tab becomes active, cursor resting on the cell we
# just edited.
$survey.Activate()
[void]$survey.Range("D5").Select()

# The "settings" tab is no longer active; its cursor was left further down
# the (otherwise empty) sheet.
[void]$settings.Range("B15").Select()

# Leave "survey" as the active/visible sheet.
$survey.Activate()
